# The sheet holds one weekly price record per row (Espinaca / Femacal de
# La Calera). A new week's record is inserted at row 477 (pushing the
# existing rows 477-511 down to 478-512), growing the sheet from 511 to
# 512 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 477; Excel shifts the old rows 477:511 down to 478:512
# and the used range/dimension grows to R512 automatically.
$ws.Rows("477:477").Insert()

# Populate the newly inserted row 477 with the new weekly record.
$ws.Range("A477").Value = 3
$ws.Range("B477").Value = "Femacal de La Calera"
$ws.Range("C477").Value = "Coquimbo"
$ws.Range("D477").Value = 45021
$ws.Range("E477").Value = 5
$ws.Range("F477").Value = 100112012
$ws.Range("G477").Value = "Espinaca"
$ws.Range("H477").Value = "Sin especificar"
$ws.Range("I477").Value = "Primera"
$ws.Range("J477").Value = 110
$ws.Range("K477").Value = 5000
$ws.Range("L477").Value = 5000
$ws.Range("M477").Value = 5000
$ws.Range("N477").Value = '$/docena de atados (3 kilos)'
$ws.Range("O477").Value = "Provincia de Quillota"
$ws.Range("P477").Value = 1667
$ws.Range("Q477").Value = 3
$ws.Range("R477").Value = "Hortaliza"
